$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$sub3 = [char]0x2083

Set-TextValue "D2" '67.099.10'
$ws.Range("E2").Value = '  -0.02%  '
Set-TextValue "D3" '3.108.87'
$ws.Range("E3").Value = '  +0.11%  '
Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue "D5" '579.07'
$ws.Range("E5").Value = '  -0.18%  '
Set-TextValue "D6" '172.88'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.69%  '
Set-TextValue "D9" '6.52'
$ws.Range("E9").Value = '  +1.24%  '
Set-TextValue "D10" '0.154'
$ws.Range("E10").Value = '  -1.08%  '
Set-TextValue "D11" '0.478'
$ws.Range("E11").Value = '  -0.69%  '
Set-TextValue "D12" '0.0000248'
Set-TextValue "D13" '36.81'
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("E14").Value = '  -1.58%  '
Set-TextValue "D15" '3.621.81'
$ws.Range("E15").Value = '  +0.24%  '
Set-TextValue "D16" '67.014.40'
$ws.Range("E16").Value = '  +0.05%  '
Set-TextValue "D17" '7.09'
$ws.Range("E17").Value = '  -1.69%  '
Set-TextValue "D18" '3.109.41'
$ws.Range("E18").Value = '  +0.30%  '
Set-TextValue "D19" '16.45'
$ws.Range("E19").Value = '  +1.03%  '
Set-TextValue "D20" '490.04'
$ws.Range("E20").Value = '  +0.88%  '
Set-TextValue "D21" '7.92'
$ws.Range("E21").Value = '  +4.57%  '
$ws.Range("E22").Value = '  -1.92%  '
Set-TextValue "D23" '83.89'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("E25").Value = '  -3.69%  '
$ws.Range("E26").Value = '  +4.59%  '
$ws.Range("E27").Value = '  -0.05%  '
Set-TextValue "D28" '7.94'
$ws.Range("E28").Value = '  -0.94%  '
$ws.Range("E29").Value = '  -1.63%  '
$ws.Range("E30").Value = '  -0.57%  '
Set-TextValue "D31" '28.35'
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("E32").Value = '  -0.60%  '
Set-TextValue "D33" "0.0${sub3}0945"
$ws.Range("E33").Value = '  -6.35%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("E36").Value = '  -1.87%  '
$ws.Range("E37").Value = '  -1.79%  '
$ws.Range("E38").Value = '  -3.97%  '
Set-TextValue "D39" '0.309'
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("E40").Value = '  +0.63%  '
Set-TextValue "D41" '8.46'
$ws.Range("E41").Value = '  -2.48%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D42" '2.802.99'
$ws.Range("E42").Value = '  -1.70%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D43" '382.81'
$ws.Range("E43").Value = '  -0.51%  '
Set-TextValue "D44" '2.58'
$ws.Range("E44").Value = '  -8.13%  '
$ws.Range("E45").Value = '  -2.64%  '
Set-TextValue "D46" '135.51'
$ws.Range("E46").Value = '  -0.30%  '
Set-TextValue "D48" '24.85'
$ws.Range("E48").Value = '  -0.83%  '
$ws.Range("E49").Value = '  -1.85%  '
$ws.Range("E50").Value = '  -1.38%  '
$ws.Range("E51").Value = '  -2.29%  '

Write-Host "Applied all changes"
